$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price record is inserted as row 97; existing rows 97-110
# (dated data, newest sighted date 45204) shift down to 98-111.
$ws.Rows("97:97").Insert()

$ws.Cells.Item(97, 1).Value = 5
$ws.Cells.Item(97, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(97, 3).Value = "Maule"
$ws.Cells.Item(97, 4).Value = 45204
$ws.Cells.Item(97, 5).Value = 7
$ws.Cells.Item(97, 6).Value = 300000000
$ws.Cells.Item(97, 7).Value = "Espárragos"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 4000
$ws.Cells.Item(97, 11).Value = 1100
$ws.Cells.Item(97, 12).Value = 1200
$ws.Cells.Item(97, 13).Value = 1150
$ws.Cells.Item(97, 14).Value = "$/kilo"
$ws.Cells.Item(97, 15).Value = "Provincia de Linares"
$ws.Cells.Item(97, 16).Value = 1150
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
